# Fix list level numbering:
# Top-level lists should be at the same paragraph level as top-level
# paragraphs (level 1 in the 1-based PowerPoint COM IndentLevel, i.e.
# lvl="0" in OOXML). Only continuation paragraphs of nested lists should
# be incremented further.

$p = $ppt.ActivePresentation

# --- Slide 1, shape 2 (Content Placeholder) ---
$s1 = $p.Slides.Item(1)
$tr1 = $s1.Shapes.Item(2).TextFrame.TextRange

# Paragraph 2: "Bullet item with inline code" : lvl 1 -> 0 (IndentLevel 2 -> 1)
$tr1.Paragraphs(2).IndentLevel = 1

# Paragraph 4: "with nested" : lvl 2 -> 1 (IndentLevel 3 -> 2)
$tr1.Paragraphs(4).IndentLevel = 2

# --- Slide 2, shape 2 (Content Placeholder) ---
$s2 = $p.Slides.Item(2)
$tr2 = $s2.Shapes.Item(2).TextFrame.TextRange

# Paragraph 2: "Nested" : lvl 1 -> 0 (IndentLevel 2 -> 1)
$tr2.Paragraphs(2).IndentLevel = 1

# --- Slide 3, shape 2 (Content Placeholder) ---
$s3 = $p.Slides.Item(3)
$tr3 = $s3.Shapes.Item(2).TextFrame.TextRange

# Paragraph 1: "A total alternative for head" : lvl 1 -> 0 (IndentLevel 2 -> 1)
$tr3.Paragraphs(1).IndentLevel = 1
